$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Status text: "In Translation" -> "Ready for handoff" ---
# (shared by Overview!E2/F2 and the "Status" column on both locale sheets)
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("C2").Value = "Ready for handoff"

# --- "Latest HO Xliff Generate Date" / de-de "Latest Handoff Datetime" ---
# "2016-10-19 11:39:00" -> "2016-10-19 11:40:14"
$wsOverview.Range("G2").Value = "2016-10-19 11:40:14"
$wsDeDe.Range("H2").Value = "2016-10-19 11:40:14"

# --- zh-cn "Latest Handoff Datetime" ---
# "2016-10-19 11:38:50" -> "2016-10-19 11:40:04"
$wsZhCn.Range("H2").Value = "2016-10-19 11:40:04"

# --- Widen the Status columns to fit "Ready for handoff" ---
# Overview columns E ("zh-cn") and F ("de-de"); zh-cn/de-de column C ("Status")
$wsOverview.Columns.Item(5).ColumnWidth = 16.33
$wsOverview.Columns.Item(6).ColumnWidth = 16.33
$wsZhCn.Columns.Item(3).ColumnWidth = 16.33
$wsDeDe.Columns.Item(3).ColumnWidth = 16.33
